$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Plain value updates (no formula involved)
$ws.Range("I12").Value = 4188377156.3100057
$ws.Range("J12").Value = 3588029419

$ws.Range("I13").Value = 1012006300.0300001
$ws.Range("J13").Value = 956934340.60000002

$ws.Range("I14").Value = -44319159.290000051
$ws.Range("J14").Value = 146269235.09999999

$ws.Range("I16").Value = -162861893.56999999
$ws.Range("J16").Value = -193292161.30000001

$ws.Range("I19").Value = -1160500000.0000002

$ws.Range("I26").Value = 1010658958.9880759
$ws.Range("J26").Value = 1031977291

# I18 and I21 gain SUM formulas (previously hard-coded values)
$ws.Range("I18").Formula = "=SUM(I12:I17)"
$ws.Range("I21").Formula = "=SUM(I18:I20)"

# Recalculate so dependent formula cells (J18, J21, I23, J23, I25, J25, J28, I29) pick up new values
$excel.CalculateFullRebuild()
$wb.Save()
